$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.651.95"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "1.596.88"
$ws.Range("E3").Value = "  +0.59%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.70"
$ws.Range("D5").ClearFormats()
$ws.Range("E6").Value = "  +1.07%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.0618"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.34%  "
$ws.Range("E9").Value = "  -0.17%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.54"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0838"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.43%  "
$ws.Range("D12").Value = "1.820.94"
$ws.Range("E12").Value = "  +0.58%  "
$ws.Range("D13").Value = "1.574.68"
$ws.Range("E13").Value = "  -1.38%  "
$ws.Range("E14").Value = "  +0.18%  "
$ws.Range("E15").Value = "  +0.24%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.47"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.54%  "
$ws.Range("D17").Value = "26.634.36"
$ws.Range("E17").Value = "  +0.12%  "
$ws.Range("E18").Value = "  +0.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "209.06"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("E20").Value = "  +0.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.96"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +3.95%  "
$ws.Range("E22").Value = "  +0.53%  "
$ws.Range("E23").Value = "  -2.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.89"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.50%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.27"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.15%  "
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.16"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.72%  "
$ws.Range("E28").Value = "  +0.69%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.29"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("E30").Value = "  +0.09%  "
$ws.Range("E31").Value = "  +0.98%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.24"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.45%  "
$ws.Range("E33").Value = "  -2.72%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.93"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.75%  "
$ws.Range("D35").Value = "1.277.55"
$ws.Range("E35").Value = "  -2.69%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.45"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.56%  "
$ws.Range("E37").Value = "  +1.10%  "
$ws.Range("E38").Value = "  -0.34%  "
$ws.Range("E39").Value = "  +2.33%  "
$ws.Range("E40").Value = "  +0.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.47"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.55%  "
$ws.Range("E42").Value = "  +1.78%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "64.60"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +3.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.786"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.72%  "
$ws.Range("E45").Value = "  +9.42%  "
$ws.Range("D46").Value = "1.733.30"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.07"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.90%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.61"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.40%  "
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.102"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +4.59%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0507"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.45%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.48"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.24%  "
